$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-63 down to 51-64
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly data record
$ws.Cells.Item(50, 1).Value = 5
$ws.Cells.Item(50, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(50, 3).Value = "Maule"
$ws.Cells.Item(50, 4).Value = 44524
$ws.Cells.Item(50, 5).Value = 7
$ws.Cells.Item(50, 6).Value = 100112026
$ws.Cells.Item(50, 7).Value = "Haba"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 200
$ws.Cells.Item(50, 11).Value = 7000
$ws.Cells.Item(50, 12).Value = 7000
$ws.Cells.Item(50, 13).Value = 7000
$ws.Cells.Item(50, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(50, 15).Value = "Región del Maule"
$ws.Cells.Item(50, 16).Value = 280
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
